$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-19 Tuesday" "2025-08-20 Wednesday"
Replace-Text "46+37=83" "91-86=5"
Replace-Text "83-79=4" "76-27=49"
Replace-Text "91-58=33" "27-16=11"
Replace-Text "82-78=4" "48+7=55"
Replace-Text "9-1=8" "99-46=53"
Replace-Text "17+9=26" "11+10=21"
Replace-Text "74-6=68" "81-65=16"
Replace-Text "55+39=94" "10+38=48"
Replace-Text "74+10=84" "91-9=82"
Replace-Text "33+46=79" "26+14=40"
Replace-Text "33+27=60" "40-6=34"
Replace-Text "3+25=28" "18+17=35"
Replace-Text "41+14=55" "65+34=99"
Replace-Text "13+30=43" "86+12=98"
Replace-Text "72+17=89" "18+54=72"
Replace-Text "58-33=25" "87+2=89"
Replace-Text "3+47=50" "12+18=30"
Replace-Text "99-95=4" "85-76=9"
Replace-Text "73-45=28" "85-68=17"
Replace-Text "17-7=10" "92-31=61"
Replace-Text "70-51=19" "38+48=86"
Replace-Text "73-11=62" "10+43=53"
Replace-Text "81-2=79" "58-51=7"
Replace-Text "70-4=66" "19+48=67"
Replace-Text "37+32=69" "29-12=17"
Replace-Text "97-18=79" "0+40=40"
Replace-Text "23+41=64" "27-16=11"
Replace-Text "42+44=86" "96-36=60"
Replace-Text "55-42=13" "5+13=18"
Replace-Text "79-3=76" "35+33=68"
Replace-Text "33-21=12" "8+54=62"
Replace-Text "74-18=56" "8+1=9"
Replace-Text "95-21=74" "50-47=3"
Replace-Text "5+3=8" "51-1=50"
Replace-Text "99-13=86" "81-29=52"
Replace-Text "31+34=65" "30+48=78"
Replace-Text "40-16=24" "1+5=6"
Replace-Text "23-11=12" "97-9=88"
Replace-Text "90-43=47" "46-5=41"
Replace-Text "92-66=26" "87-36=51"
Replace-Text "29-15=14" "42+23=65"
Replace-Text "88-74=14" "38+6=44"
Replace-Text "85-14=71" "59-0=59"
Replace-Text "14-6=8" "77-35=42"
Replace-Text "28+0=28" "35-24=11"
Replace-Text "40+14=54" "47+20=67"
Replace-Text "56+16=72" "70-33=37"
Replace-Text "10+74=84" "48+47=95"
Replace-Text "24+39=63" "60+20=80"
Replace-Text "75+12=87" "79-36=43"
Replace-Text "3+67=70" "65-13=52"
Replace-Text "24+18=42" "48+39=87"
Replace-Text "51+1=52" "46-2=44"
Replace-Text "29+24=53" "52-15=37"
Replace-Text "75-15=60" "45+49=94"
Replace-Text "19+9=28" "25-10=15"
Replace-Text "69-19=50" "85-80=5"
Replace-Text "95-74=21" "83-82=1"
Replace-Text "54+29=83" "27-19=8"
Replace-Text "5+17=22" "93+0=93"
Replace-Text "5+75=80" "95-67=28"
Replace-Text "90-22=68" "93+0=93"
Replace-Text "53-20=33" "40+23=63"
Replace-Text "35+2=37" "43+30=73"
Replace-Text "99-16=83" "94-90=4"
Replace-Text "42+57=99" "5+51=56"
Replace-Text "44+32=76" "16+59=75"
Replace-Text "99-24=75" "38-21=17"
Replace-Text "77+17=94" "73-39=34"
Replace-Text "44+4=48" "43+16=59"
Replace-Text "15+75=90" "35+55=90"
Replace-Text "38+26=64" "86-60=26"
Replace-Text "92-49=43" "13+71=84"
Replace-Text "88-7=81" "78+18=96"
Replace-Text "29+15=44" "65-14=51"
Replace-Text "76-9=67" "75+23=98"
Replace-Text "72-50=22" "80-79=1"
Replace-Text "58+40=98" "93-51=42"
Replace-Text "46-38=8" "88-28=60"
Replace-Text "76-50=26" "63-32=31"
Replace-Text "70-29=41" "38+51=89"
Replace-Text "15+76=91" "34+15=49"
Replace-Text "62-21=41" "83-57=26"
Replace-Text "69+3=72" "74-13=61"
Replace-Text "71+7=78" "13+81=94"
Replace-Text "33+38=71" "57+26=83"
Replace-Text "91-48=43" "73-5=68"
Replace-Text "37+10=47" "7+89=96"
Replace-Text "10+35=45" "39-11=28"
Replace-Text "44+12=56" "57-4=53"
Replace-Text "94-18=76" "42+35=77"
Replace-Text "29+2=31" "87-75=12"
Replace-Text "11+23=34" "10+3=13"
Replace-Text "8+83=91" "59-34=25"
Replace-Text "37-7=30" "60+29=89"
Replace-Text "94-34=60" "94-36=58"
Replace-Text "50-44=6" "38+47=85"
Replace-Text "65-21=44" "88+11=99"
Replace-Text "4+8=12" "0+84=84"
Replace-Text "12+80=92" "74-59=15"
